$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge duplicate sale rows: the duplicate entries in rows 15-18 are no
# longer needed (they were exact copies of row 14), so clear their
# contents while keeping the existing cell formatting.
$ws.Range("A15:I18").ClearContents()

# Update the active selection to match the new focus point.
$ws.Range("B13").Select()
